$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The config sheet had three settings rows that are no longer used by the
# automation (the log file was never actually being updated, so the
# screenshots/log/max-tickers entries were dropped when logging was
# reworked):
#   Row 3: ScreenshotsFolder / Screenshots
#   Row 4: LogFilePath / Logs\run.log
#   Row 6: MaxTickersToProcess / 20
#
# Remove the existing hyperlink objects first -- they are anchored to fixed
# cells and won't otherwise follow the rows they decorate when those rows
# shift up after the deletions below.
$ws.Hyperlinks.Delete()

# Delete bottom-most row first so the remaining row numbers used below stay
# valid.
$ws.Rows("6:6").Delete()
$ws.Rows("3:4").Delete()

# Re-create the three hyperlinks against their new (shifted-up) cells, in
# the same order they originally appeared so the relationship ids line up
# the same way (rId1 -> most-active link, rId2 -> Yahoo Finance homepage,
# rId3 -> mailto recipient).
#
# B5's literal cell text (the ExtractDataFromUrl JSON array) is not the
# same as this hyperlink's display text, so the saved file keeps a
# "display" override alongside the link. Pass the display text through
# TextToDisplay (which also stamps it onto the cell), then restore the
# cell's real text afterwards.
$b5Text = $ws.Range("B5").Value2
$ws.Hyperlinks.Add($ws.Range("B5"), "https://ca.finance.yahoo.com/markets/stocks/most-active/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://ca.finance.yahoo.com/markets/stocks/most-active/")
$ws.Range("B5").Value2 = $b5Text

$ws.Hyperlinks.Add($ws.Range("B4"), "https://finance.yahoo.com/")
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:matan10cohen@gmail.com")

# Re-adding hyperlinks resets the formatting on those cells -- restore the
# workbook's "Hyperlink" cell style so the cells look exactly as before.
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"

# Match the saved selection (whole row 3 selected) reflected in the file.
$ws.Rows("3:3").Select()
